$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 72.55556
$ws.Range("I5").Value = 75.375
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 75.375
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 39.625
$ws.Range("N5").Value = -280

$ws.Range("H17").Value = 26500
$ws.Range("J17").Value = 26500
$ws.Range("L17").Value = 79500
$ws.Range("N17").Value = -79836

$ws.Range("H18").Value = 1398.1666
$ws.Range("I18").Value = 878
$ws.Range("J18").Value = 3999
$ws.Range("K18").Value = 878
$ws.Range("L18").Value = 3999
$ws.Range("M18").Value = -594
$ws.Range("N18").Value = -4567

$ws.Range("H29").Value = 5829.5713
$ws.Range("I29").Value = 4933.3335
$ws.Range("J29").Value = 6501.75
$ws.Range("K29").Value = 14800.0005
$ws.Range("L29").Value = 19505.25
$ws.Range("M29").Value = -14519.0005
$ws.Range("N29").Value = -20067.25

$ws.Range("H46").Value = 4100
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15238

$ws.Range("H52").Value = 795.75
$ws.Range("I52").Value = 795.75
$ws.Range("K52").Value = 2387.25
$ws.Range("M52").Value = -2227.25

$ws.Range("H60").Value = 4100
$ws.Range("J60").Value = 5000
$ws.Range("L60").Value = 15000
$ws.Range("N60").Value = -15968

$ws.Range("H61").Value = 2807.5
$ws.Range("I61").Value = 2410
$ws.Range("K61").Value = 7230
$ws.Range("M61").Value = -7058

$ws.Range("H112").Value = 2589.125
$ws.Range("J112").Value = 2589.125
$ws.Range("L112").Value = 7767.375
$ws.Range("N112").Value = -9983.375

$ws.Range("H133").Value = 69895.55
$ws.Range("J133").Value = 69895.55
$ws.Range("L133").Value = 69895.55
$ws.Range("N133").Value = -80015.55

$ws.Range("H134").Value = 67463.664
$ws.Range("J134").Value = 67463.664
$ws.Range("L134").Value = 67463.664
$ws.Range("N134").Value = -77603.664

$ws.Range("H136").Value = 68899
$ws.Range("J136").Value = 68899
$ws.Range("L136").Value = 68899
$ws.Range("N136").Value = -79099

$ws.Range("H138").Value = 3647.724
$ws.Range("J138").Value = 4102.804
$ws.Range("L138").Value = 12308.412
$ws.Range("N138").Value = -22588.412

$ws.Range("H139").Value = 48403.5
$ws.Range("J139").Value = 48403.5
$ws.Range("L139").Value = 48403.5
$ws.Range("N139").Value = -58683.5

$ws.Range("H140").Value = 69054.5
$ws.Range("J140").Value = 69054.5
$ws.Range("L140").Value = 69054.5
$ws.Range("N140").Value = -79414.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 687.5
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 1075
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 1075
$ws.Range("M4").Value = -184
$ws.Range("N4").Value = -1307

$ws.Range("H32").Value = 5376.176
$ws.Range("I32").Value = 4373.4043
$ws.Range("K32").Value = 4373.4043
$ws.Range("M32").Value = -4086.4043

$ws.Range("H61").Value = 3126.6458
$ws.Range("J61").Value = 7445.4
$ws.Range("L61").Value = 7445.4
$ws.Range("N61").Value = -7869.4

$ws.Range("H132").Value = 3637.25
$ws.Range("I132").Value = 2239.3635
$ws.Range("K132").Value = 6718.0905
$ws.Range("M132").Value = -4188.0905

$ws.Range("H136").Value = 3126.6458
$ws.Range("J136").Value = 7445.4
$ws.Range("L136").Value = 22336.2
$ws.Range("N136").Value = -27436.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 260
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -646

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 10174.25
$ws.Range("I22").Value = 5348.5
$ws.Range("J22").Value = 15000
$ws.Range("K22").Value = 5348.5
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = -4998.5
$ws.Range("N22").Value = -15700

$ws.Range("H31").Value = 33128.83
$ws.Range("I31").Value = 1208.4445
$ws.Range("J31").Value = 66926.88
$ws.Range("K31").Value = 1208.4445
$ws.Range("L31").Value = 66926.88
$ws.Range("M31").Value = -913.4445000000001
$ws.Range("N31").Value = -67516.88

$ws.Range("H34").Value = 33128.83
$ws.Range("I34").Value = 1208.4445
$ws.Range("J34").Value = 66926.88
$ws.Range("K34").Value = 1208.4445
$ws.Range("L34").Value = 66926.88
$ws.Range("M34").Value = -1006.4445
$ws.Range("N34").Value = -67330.88

$ws.Range("H58").Value = 2989.5454
$ws.Range("I58").Value = 1531.1666
$ws.Range("K58").Value = 1531.1666
$ws.Range("M58").Value = -1328.1666

$ws.Range("H136").Value = 2989.5454
$ws.Range("I136").Value = 1531.1666
$ws.Range("K136").Value = 4593.4998
$ws.Range("M136").Value = -2043.4998

$ws.Range("H138").Value = 47570.145
$ws.Range("J138").Value = 47570.145
$ws.Range("L138").Value = 47570.145
$ws.Range("N138").Value = -57850.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5657
$ws.Range("I56").Value = 5657
$ws.Range("K56").Value = 5657
$ws.Range("M56").Value = -5127

$ws.Range("H107").Value = 473838.44
$ws.Range("J107").Value = 1250468.9
$ws.Range("L107").Value = 3751406.7
$ws.Range("N107").Value = -3755246.7

$ws.Range("H124").Value = 5030.125
$ws.Range("J124").Value = 5748.2
$ws.Range("L124").Value = 17244.6
$ws.Range("N124").Value = -27064.6

$ws.Range("H128").Value = 178233.25
$ws.Range("I128").Value = 178233.25
$ws.Range("K128").Value = 534699.75
$ws.Range("M128").Value = -529719.75

$ws.Range("H132").Value = 4448.0625
$ws.Range("I132").Value = 3213.889
$ws.Range("K132").Value = 28925.001
$ws.Range("M132").Value = -26395.001

$ws.Range("H133").Value = 17550126
$ws.Range("I133").Value = 1243.5714
$ws.Range("J133").Value = 27786974
$ws.Range("K133").Value = 3730.7142
$ws.Range("L133").Value = 83360922
$ws.Range("M133").Value = 1329.2858
$ws.Range("N133").Value = -83371042

$ws.Range("H139").Value = 4876.25
$ws.Range("I139").Value = 2500
$ws.Range("K139").Value = 7500
$ws.Range("M139").Value = -2360

$ws.Range("H141").Value = 11036.25
$ws.Range("I141").Value = 6939.727
$ws.Range("J141").Value = 14502.538
$ws.Range("K141").Value = 20819.181
$ws.Range("L141").Value = 43507.614
$ws.Range("M141").Value = -15639.181
$ws.Range("N141").Value = -53867.614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 446.86957
$ws.Range("I2").Value = 50.4
$ws.Range("J2").Value = 1190.25
$ws.Range("K2").Value = 50.4
$ws.Range("L2").Value = 1190.25
$ws.Range("M2").Value = 62.6
$ws.Range("N2").Value = -1416.25

$ws.Range("H126").Value = 3391.276
$ws.Range("I126").Value = 2093
$ws.Range("J126").Value = 4782.2856
$ws.Range("K126").Value = 6279
$ws.Range("L126").Value = 14346.8568
$ws.Range("M126").Value = -3809
$ws.Range("N126").Value = -19286.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3428.3684
$ws.Range("I22").Value = 1045.5
$ws.Range("J22").Value = 10100.4
$ws.Range("K22").Value = 1045.5
$ws.Range("L22").Value = 10100.4
$ws.Range("M22").Value = -750.5
$ws.Range("N22").Value = -10690.4

$ws.Range("H27").Value = 3428.3684
$ws.Range("I27").Value = 1045.5
$ws.Range("J27").Value = 10100.4
$ws.Range("K27").Value = 1045.5
$ws.Range("L27").Value = 10100.4
$ws.Range("M27").Value = -938.5
$ws.Range("N27").Value = -10314.4

$ws.Range("H46").Value = 4900.125
$ws.Range("I46").Value = 1399.6666
$ws.Range("K46").Value = 1399.6666
$ws.Range("M46").Value = -1211.6666

$ws.Range("H93").Value = 2049.5
$ws.Range("I93").Value = 2145.125
$ws.Range("K93").Value = 2145.125
$ws.Range("M93").Value = -897.125

$ws.Range("H132").Value = 5997.375
$ws.Range("I132").Value = 4583.85
$ws.Range("K132").Value = 13751.55
$ws.Range("M132").Value = -11221.55

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2256.4
$ws.Range("I122").Value = 1665.303
$ws.Range("K122").Value = 4995.909000000001
$ws.Range("M122").Value = -2545.909000000001
